# Update workbook per commit: "Add message box when nhap, xuat, edit,
# search names, export file Excel" — underlying data change: the unit
# price on the "nhap-linhkien" (import) sheet is corrected, and a new
# row is recorded on the "xuat-linhkien" (export) sheet.

$wb = $excel.ActiveWorkbook

$wsNhap = $wb.Worksheets.Item("nhap-linhkien")
$wsXuat = $wb.Worksheets.Item("xuat-linhkien")

# --- nhap-linhkien: the part referenced on row 2 was renamed/corrected,
# and its unit price (Đơn Giá) was updated.
$wsNhap.Range("A2").Value = "PCB - RF_1pha_HT32F5_Si4432_E19_E49_ESRF_V1.3- 03022021 Size 45x65mm"
$wsNhap.Range("B2").Value = "LF_1pha_HT32F5 V1.3_03022021"
$wsNhap.Range("D2").Value = "MODULE RF 1P HOLTEK"
$wsNhap.Range("H2").Value = 266704

# --- xuat-linhkien: append a new export record on row 2 ---
$wsXuat.Range("A2").Value = "PCB - RF_1pha_HT32F5_Si4432_E19_E49_ESRF_V1.3- 03022021 Size 45x65mm"
$wsXuat.Range("B2").Value = "LF_1pha_HT32F5 V1.3_03022021"

# Sổ Hợp Đồng (contract number) is blank for this record.
$wsXuat.Range("C2").Value = ""

$wsXuat.Range("D2").Value = "MODULE RF 1P HOTEK"
$wsXuat.Range("E2").Value = "CTY TNHH YEAR2000"

# Ngày Nhập (date) must stay a plain text value like "2021-06-30" on the
# other sheet, not get auto-converted into a date serial number.
$wsXuat.Range("F2").NumberFormat = "@"
$wsXuat.Range("F2").Value = "2021-07-01"
$wsXuat.Range("F2").Style = "Normal"

$wsXuat.Range("G2").Value = "Cái"
$wsXuat.Range("H2").Value = 100
$wsXuat.Range("I2").Value = 0
$wsXuat.Range("J2").Value = 0
